$d = $word.ActiveDocument

# --- 1. Remove the "Meta description" paragraph that follows the title heading ---
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# --- 2. Replace the final (image-prompt) paragraph with two paragraphs:
#        a bold "Play Football Glory for Free - Slot Review" paragraph, and
#        an italic meta-description paragraph. ---
$count = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($count)
$r = $pLast.Range

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Football Glory for Free - Slot Review</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover Football Glory, a football-themed slot with captivating graphics and unique gameplay. Play for free and explore special symbols like fixed Wilds and Cup symbol.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# InsertXML leaves a stray empty paragraph at the very end of the body (since the
# inserted content ends with a paragraph mark of its own); remove that trailing
# empty paragraph by deleting from the end of the previous paragraph through it.
$countNow = $d.Paragraphs.Count
$pPrev = $d.Paragraphs.Item($countNow - 1)
$pExtra = $d.Paragraphs.Item($countNow)
$rDel = $d.Range($pPrev.Range.End - 1, $pExtra.Range.End)
$rDel.Delete()
